{"js": "// Change the year in the astromap link: \".../GaNight/2018/\" -> \".../GaNight/2022/\"\n// The credit/link paragraph's text is located by a full-text search, then its\n// entire run content is replaced (formatting collapses to a single plain run,\n// matching how Word behaves when the whole paragraph text is retyped/pasted).\nconst body = context.document.body;\n\nconst oldText =\n  \"Les cartes figurant dans ce document ont \u00e9t\u00e9 \u00e9tablies par Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/).\";\nconst newText =\n  \"Les cartes figurant dans ce document ont \u00e9t\u00e9 \u00e9tablies par Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\";\n\nconst searchResults = body.search(oldText, { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the astromap credit paragraph text.\");\n}\n\nconst target = searchResults.items[0];\n\n// Escape the few XML-sensitive characters that can appear in the text.\nfunction xmlEscape(text) {\n  return text\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p><w:r/><w:r><w:t>\" +\n  xmlEscape(newText) +\n  \"</w:t></w:r></w:p></w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that contains the astromap credit/link text.\n$target = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"*Les cartes figurant*\" -and $t -like \"*GaNight/2018*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $r = $target.Range\n    # Exclude the trailing paragraph mark from the replaced range.\n    $body = $d.Range($r.Start, $r.End - 1)\n\n    $newText = \"Les cartes figurant dans ce document ont \u00e9t\u00e9 \u00e9tablies par Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\"\n\n    $newXml = '<?xml version=\"1.0\" standalone=\"yes\"?><?mso-application progid=\"Word.Document\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r/><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n    $body.InsertXML($newXml)\n}\n"}
